# Mise à jour du classement - 03.04.2025 à 17:00
$wb = $excel.ActiveWorkbook

# --- leaderboard2 ---
$ws2 = $wb.Worksheets.Item("leaderboard2")
$ws2.Range("D4").Value = 555
$ws2.Range("B13").Value = "Dernière update le 03.04.25 à 17:00"

# --- leaderboard3 ---
$ws3 = $wb.Worksheets.Item("leaderboard3")
$ws3.Range("B13").Value = "Dernière update le 03.04.25 à 17:00"

# --- leaderboard4 ---
$ws4 = $wb.Worksheets.Item("leaderboard4")
$ws4.Range("C3").Value = "BKZRackham"
$ws4.Range("D3").Value = 21
$ws4.Range("C4").Value = "ArtyumsM"
$ws4.Range("D4").Value = 20
$ws4.Range("B13").Value = "Dernière update le 03.04.25 à 17:00"
